$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the missing survey-3 data row (row 4: B,C,D)
$ws.Range("B4").Value = 9
$ws.Range("C4").Value = 22
$ws.Range("D4").Value = 0

# Update the selection to match the newly entered range
$ws.Range("B4:D4").Select()
